$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: change highlight color from red to green on the paragraph
# "3. Mostrar en la web el balance de BNB que hay en la wallet del
# usuario (0,4)" (including its paragraph mark, so the pPr/rPr is
# updated too).
# ---------------------------------------------------------------------
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Mostrar en la web el balance de BNB") {
        $target1 = $para
        break
    }
}
$target1.Range.Font.HighlightColorIndex = 4

# ---------------------------------------------------------------------
# Edit 2: change "Total: 0,9 puntos" to "Total: 1,3 puntos" while
# preserving the original run structure (5 runs) and their rsid
# attributes. A plain text replace would coalesce the runs, so the
# paragraph content (excluding its trailing paragraph mark) is
# replaced via InsertXML with the exact desired run structure.
# ---------------------------------------------------------------------
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^Total: 0,9 puntos") {
        $target2 = $para
        break
    }
}
$full = $target2.Range
$body = $d.Range($full.Start, $full.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Total: </w:t></w:r><w:r w:rsidR="0000385A"><w:t>1</w:t></w:r><w:r w:rsidR="003B0951"><w:t>,</w:t></w:r><w:r w:rsidR="00186B29"><w:t>3</w:t></w:r><w:r w:rsidR="0000385A"><w:t xml:space="preserve"> puntos</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$body.InsertXML($xml)
